$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Edit Deal test cases (rows 8-10) ----
$ws.Cells.Item(8, 1).Value = "EditDeal_TC001"
$ws.Cells.Item(8, 2).Value = "John Tucker"
$ws.Cells.Item(8, 3).Value = "ONE"
$ws.Cells.Item(8, 4).Value = "Deal shared successfully"

$ws.Cells.Item(9, 1).Value = "EditDeal_TC002"
$ws.Cells.Item(9, 2).Value = "NA"
$ws.Cells.Item(9, 3).Value = "ALL"
$ws.Cells.Item(9, 4).Value = "Deal shared successfully"

$ws.Cells.Item(10, 1).Value = "EditDeal_TC003"
$ws.Cells.Item(10, 2).Value = "Stan Koster Andersons"
$ws.Cells.Item(10, 3).Value = "ONE"
$ws.Cells.Item(10, 4).Value = "Deal shared successfully"

# ---- ReShare Deal test cases, entered as Share/ReShare pairs (rows 11-16) ----
$ws.Cells.Item(11, 1).Value = "ReShareDeal_TC001"
$ws.Cells.Item(12, 1).Value = "ReShareDeal_TC001(2)"
$ws.Cells.Item(11, 2).Value = "John Tucker"
$ws.Cells.Item(12, 2).Value = "John Tucker"
$ws.Cells.Item(11, 3).Value = "ONE"
$ws.Cells.Item(12, 3).Value = "RESHARE"
$ws.Cells.Item(11, 4).Value = "Deal shared successfully"
$ws.Cells.Item(12, 4).Value = "Deal shared successfully"

$ws.Cells.Item(13, 1).Value = "ReShareDeal_TC002"
$ws.Cells.Item(14, 1).Value = "ReShareDeal_TC002(2)"
$ws.Cells.Item(13, 2).Value = "Stan Koster Andersons"
$ws.Cells.Item(14, 2).Value = "Stan Koster Andersons"
$ws.Cells.Item(13, 3).Value = "ONE"
$ws.Cells.Item(14, 3).Value = "RESHARE"
$ws.Cells.Item(13, 4).Value = "Deal shared successfully"
$ws.Cells.Item(14, 4).Value = "Deal shared successfully"

$ws.Cells.Item(15, 1).Value = "ReShareDeal_TC003"
$ws.Cells.Item(16, 1).Value = "ReShareDeal_TC003(2)"
$ws.Cells.Item(15, 2).Value = "Karthikeyan"
$ws.Cells.Item(16, 2).Value = "Karthikeyan"
$ws.Cells.Item(15, 3).Value = "ONE"
$ws.Cells.Item(16, 3).Value = "RESHARE"
$ws.Cells.Item(15, 4).Value = "Deal shared successfully"
$ws.Cells.Item(16, 4).Value = "Deal shared successfully"

# Match formatting of existing data rows (vertical-center style) for column A
$ws.Range("A8:A16").VerticalAlignment = -4108

$ws.Columns.Item(1).EntireColumn.AutoFit() | Out-Null
